$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

$rows1 = @(
    @("BRVM - SERVICES PUBLICS", 0, 8, 3360.36, 106.37, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI", 0, 4, 2710, 655, "🟡 Observer", "➖ Neutre"),
    @("AIR LIQUIDE CI", 0, 4, 2625, 670, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AUTRES SECTEURS", 0, 4, 2371.77, 605.65, "🟡 Observer", "➖ Neutre"),
    @("BRVM - DISTRIBUTION", 0, 4, 1861.76, 478.74, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AGRICULTURE", 0, 4, 1491.94, 376.98, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TRANSPORT", 0, 4, 1442.94, 362.26, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 4, 618.96, 161.4, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRESTIGE", 0, 4, 548.47, 137.58, "🟡 Observer", "➖ Neutre"),
    @("BRVM - FINANCES", 0, 4, 542.03, 135.61, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES FINANCIERS", 0, 4, 532.71, 133.28, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIELS", 0, 4, 512.47, 130, "🟡 Observer", "➖ Neutre"),
    @("BRVM - ENERGIE", 0, 4, 438.34, 110.17, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TELECOMMUNICATIONS", 0, 4, 383.03, 95.99, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                   (**)", 0, 1, 213.25, 213.25, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                        (**)", 0, 1, 212.57, 212.57, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                           (**)", 0, 1, 199.2, 199.2, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                      (**)", 0, 1, 198.63, 198.63, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE             (**)", 0, 1, 192.99, 192.99, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE              (**)", 0, 1, 192.06, 192.06, "🟡 Observer", "➖ Neutre"),
    @("TRACTAFRIC MOTORS CI (PRSC)", 4, 0, 29.5, 7.37, "🟢 Achat", "✅ Renforcer"),
    @("CFAO MOTORS CI (CFAC)", 2, 0, 14.81, 7.35, "🟡 Observer", "➖ Neutre"),
    @("SICABLE CI (CABC)", 2, 0, 14.78, 7.45, "🟡 Observer", "➖ Neutre"),
    @("FILTISAC CI (FTSC)", 3, 1, 14.7, -7.46, "🟢 Achat", "✅ Renforcer"),
    @("SETAO CI (STAC)", 2, 1, 11.68, 7.27, "🟡 Observer", "👀 À surveiller"),
    @("BERNABE CI (BNBC)", 2, 1, 10.85, -3.16, "🟡 Observer", "👀 À surveiller"),
    @("SAFCA CI (SAFC)", 1, 0, 7.34, 7.34, "🟡 Observer", "➖ Neutre"),
    @("UNIWAX CI (UNXC)", 1, 0, 7.23, 7.23, "🟡 Observer", "➖ Neutre"),
    @("SUCRIVOIRE (SCRC)", 1, 0, 6.78, 6.78, "🟡 Observer", "➖ Neutre"),
    @("ONATEL BF (ONTBF)", 1, 1, 5.41, 7.31, "🟡 Observer", "👀 À surveiller"),
    @("SAPH CI (SPHC)", 1, 1, 0.64, -6.81, "🟡 Observer", "👀 À surveiller"),
    @("TOTAL", 0, 5, 0, 0, "🟡 Observer", "➖ Neutre"),
    @("BANK OF AFRICA ML (BOAM)", 0, 1, -1.15, -1.15, "🟡 Observer", "➖ Neutre"),
    @("SOGB CI (SOGC)", 0, 1, -1.18, -1.18, "🟡 Observer", "➖ Neutre"),
    @("SOCIETE GENERALE COTE D'IVOIRE (SGBC)", 0, 1, -1.85, -1.85, "🟡 Observer", "➖ Neutre"),
    @("TOTALENERGIES MARKETING CI (TTLC)", 0, 1, -2, -2, "🟡 Observer", "➖ Neutre"),
    @("LOTERIE NATIONALE DU BENIN (LNBB)", 0, 1, -2.13, -2.13, "🟡 Observer", "➖ Neutre"),
    @("ECOBANK COTE D''IVOIRE (ECOC)", 0, 1, -2.23, -2.23, "🟡 Observer", "➖ Neutre"),
    @("SOLIBRA CI (SLBC)", 0, 1, -2.75, -2.75, "🟡 Observer", "➖ Neutre"),
    @("ORANGE COTE D'IVOIRE (ORAC)", 0, 1, -3.34, -3.34, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI (NEIC)", 0, 1, -3.65, -3.65, "🟡 Observer", "➖ Neutre"),
    @("NSIA BANQUE COTE D'IVOIRE (NSBC)", 0, 2, -3.84, -1.83, "🟡 Observer", "➖ Neutre"),
    @("AIR LIQUIDE CI (SIVC)", 0, 1, -7.35, -7.35, "🟡 Observer", "➖ Neutre"),
    @("ORAGROUP TOGO (ORGT)", 0, 1, -7.38, -7.38, "🟡 Observer", "➖ Neutre"),
)

$r = 2
foreach ($row in $rows1) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

$rows2 = @(
    @("BRVM - SERVICES PUBLICS", 8713819.89),
    @("NEI-CEDA CI", 365106.34),
    @("AIR LIQUIDE CI", 326476.25),
    @("BRVM - AUTRES SECTEURS", 230154.57),
    @("BRVM - DISTRIBUTION", 101997.33),
    @("BRVM - AGRICULTURE", 49909.48),
    @("BRVM - TRANSPORT", 44957.17),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 4104.16),
    @("BRVM-PRESTIGE", 3061.2),
    @("BRVM - FINANCES", 2976.23),
)

$r = 2
foreach ($row in $rows2) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}
